$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the header row ("DIV 1", "DIV 2", "DIV 3") - all remaining rows shift up by one.
$ws.Rows.Item(1).Delete()

# Select the new first row (mirrors the resulting selection after a row delete).
$ws.Rows.Item(1).Select()
